# fall 23 week 5 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("F2").Value = 10.2
$ws.Range("G2").Value = 12

$ws.Range("D3").Value = 10.61

$ws.Range("C4").Value = 9.390000000000001
$ws.Range("E4").Value = 10.21
$ws.Range("F4").Value = 10.25

$ws.Range("D5").Value = 9.789999999999999
$ws.Range("G5").Value = 9.24

$ws.Range("B6").Value = 9.800000000000001
$ws.Range("D6").Value = 9.75
$ws.Range("H6").Value = 11.24
$ws.Range("I6").Value = 9

$ws.Range("B7").Value = 8
$ws.Range("E7").Value = 10.76
$ws.Range("H7").Value = 9.6

$ws.Range("F8").Value = 8.76
$ws.Range("G8").Value = 10.4

$ws.Range("F9").Value = 11
